$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the explicit "justify" (both) alignment override from the
#    eight numbered "List Paragraph" items in the "Έχοντας υπόψη" list,
#    so they fall back to the document default (left) alignment, i.e.
#    the <w:jc w:val="both"/> element disappears from their <w:pPr>.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "List Paragraph") {
        $p.Alignment = 0
    }
}

# ---------------------------------------------------------------------
# 2. Disable the city field: "... ${school} ${country} (${city}) από ..."
#    becomes "... ${school} ${country}, από ...".
#    Edit the run that holds the literal "city" text in place (so it
#    keeps living as its own run rather than merging into its
#    neighbours), then delete the now-redundant " (${" / "})" runs
#    around it.
# ---------------------------------------------------------------------
$whole = $d.Content
$whole.Find.Execute(' (${city})', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $whole.Start
$end = $whole.End

# Drop the trailing "})" first (rightmost edit, leaves earlier offsets intact).
$suffix = $d.Range($end - 2, $end)
$suffix.Text = ""

# Turn the "city" run itself into the replacement comma.
$cityRun = $d.Range($start + 4, $start + 8)
$cityRun.Text = ","

# Drop the leading " (${" run.
$prefix = $d.Range($start, $start + 4)
$prefix.Text = ""

# ---------------------------------------------------------------------
# 3. Disable the class/title suffix: "... της ${class} τάξης." becomes
#    "... της ${class}.".
# ---------------------------------------------------------------------
$tail = $d.Content
$tail.Find.Execute('} τάξης.', $false, $false, $false, $false, $false, $true, 1, $false, '}.', 2) | Out-Null
